# Updates the cryptos list data (price/volume columns, and a couple of
# coin-name/link/value swaps) to match the latest scrape.
# Numeric-looking price strings are forced to stay as text (matching the
# original inlineStr cell type) by temporarily applying a text number
# format and then clearing the format again so no stray style is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.970.59"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").Value = "1.726.91"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  -0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "218.57"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +1.62%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.524"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.10%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "24.19"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +14.26%  "
$ws.Range("E9").Value = "  +3.69%  "
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "1.970.69"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").Value = "1.722.76"
$ws.Range("E13").Value = "  +2.81%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.566"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +6.27%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "67.90"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "27.930.85"
$ws.Range("E17").Value = "  +3.33%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "243.69"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("E19").Value = "  +2.61%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.90"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +4.52%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.78"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +4.92%  "
$ws.Range("E24").Value = "  +0.48%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "149.58"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("E26").Value = "  +4.34%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "16.82"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("D34").Value = "1.485.06"
$ws.Range("E34").Value = "  -3.68%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.67"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.961"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +4.76%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.612"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("E40").Value = "  +0.54%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "71.48"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +5.49%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.84"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").Value = "1.875.05"
$ws.Range("E45").Value = "  +3.17%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.792"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +13.48%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "91.79"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  +1.91%  "
